$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-301)
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03).
$ws.Range("C2:C301").Value = 45172
